$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A24").Value = 41369
$ws.Range("A24").NumberFormat = "m/d/yy"

$ws.Range("B24").Value = 0.0625
$ws.Range("B24").NumberFormat = "h:mm"

[void]$ws.Range("B25").Select()
